$d = $word.ActiveDocument

$d.Content.Find.Execute("2026-01-18 Sunday", $true, $false, $false, $false, $false, $true, 1, $false, "2026-01-19 Monday", 2) | Out-Null
$d.Content.Find.Execute("838÷8=104, 6", $true, $false, $false, $false, $false, $true, 1, $false, "304÷7=43, 3", 2) | Out-Null
$d.Content.Find.Execute("244÷6=40, 4", $true, $false, $false, $false, $false, $true, 1, $false, "113÷6=18, 5", 2) | Out-Null
$d.Content.Find.Execute("504÷8=63, 0", $true, $false, $false, $false, $false, $true, 1, $false, "263÷8=32, 7", 2) | Out-Null
$d.Content.Find.Execute("695÷7=99, 2", $true, $false, $false, $false, $false, $true, 1, $false, "815÷4=203, 3", 2) | Out-Null
$d.Content.Find.Execute("507÷9=56, 3", $true, $false, $false, $false, $false, $true, 1, $false, "311÷3=103, 2", 2) | Out-Null
$d.Content.Find.Execute("131÷5=26, 1", $true, $false, $false, $false, $false, $true, 1, $false, "481÷3=160, 1", 2) | Out-Null
$d.Content.Find.Execute("183÷8=22, 7", $true, $false, $false, $false, $false, $true, 1, $false, "943÷4=235, 3", 2) | Out-Null
$d.Content.Find.Execute("337÷3=112, 1", $true, $false, $false, $false, $false, $true, 1, $false, "585÷7=83, 4", 2) | Out-Null
$d.Content.Find.Execute("828÷9=92, 0", $true, $false, $false, $false, $false, $true, 1, $false, "463÷4=115, 3", 2) | Out-Null
$d.Content.Find.Execute("230÷8=28, 6", $true, $false, $false, $false, $false, $true, 1, $false, "906÷6=151, 0", 2) | Out-Null
$d.Content.Find.Execute("978÷3=326, 0", $true, $false, $false, $false, $false, $true, 1, $false, "495÷5=99, 0", 2) | Out-Null
$d.Content.Find.Execute("395÷6=65, 5", $true, $false, $false, $false, $false, $true, 1, $false, "176÷9=19, 5", 2) | Out-Null
$d.Content.Find.Execute("441÷5=88, 1", $true, $false, $false, $false, $false, $true, 1, $false, "909÷5=181, 4", 2) | Out-Null
$d.Content.Find.Execute("381÷5=76, 1", $true, $false, $false, $false, $false, $true, 1, $false, "125÷4=31, 1", 2) | Out-Null
$d.Content.Find.Execute("301÷9=33, 4", $true, $false, $false, $false, $false, $true, 1, $false, "980÷2=490, 0", 2) | Out-Null
$d.Content.Find.Execute("168÷7=24, 0", $true, $false, $false, $false, $false, $true, 1, $false, "668÷5=133, 3", 2) | Out-Null
$d.Content.Find.Execute("180÷3=60, 0", $true, $false, $false, $false, $false, $true, 1, $false, "821÷7=117, 2", 2) | Out-Null
$d.Content.Find.Execute("514÷3=171, 1", $true, $false, $false, $false, $false, $true, 1, $false, "455÷8=56, 7", 2) | Out-Null
$d.Content.Find.Execute("506÷2=253, 0", $true, $false, $false, $false, $false, $true, 1, $false, "422÷7=60, 2", 2) | Out-Null
$d.Content.Find.Execute("510÷4=127, 2", $true, $false, $false, $false, $false, $true, 1, $false, "892÷4=223, 0", 2) | Out-Null
$d.Content.Find.Execute("574÷6=95, 4", $true, $false, $false, $false, $false, $true, 1, $false, "224÷7=32, 0", 2) | Out-Null
$d.Content.Find.Execute("408÷6=68, 0", $true, $false, $false, $false, $false, $true, 1, $false, "716÷4=179, 0", 2) | Out-Null
$d.Content.Find.Execute("230÷3=76, 2", $true, $false, $false, $false, $false, $true, 1, $false, "704÷6=117, 2", 2) | Out-Null
$d.Content.Find.Execute("504÷9=56, 0", $true, $false, $false, $false, $false, $true, 1, $false, "568÷5=113, 3", 2) | Out-Null
$d.Content.Find.Execute("513÷8=64, 1", $true, $false, $false, $false, $false, $true, 1, $false, "320÷8=40, 0", 2) | Out-Null
